$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Part 1: Title - merge "Book 3 - 5" / "0 " / "- " / "On The Olive Mount"
#         into a single run "Book 3 - 50 - On The Olive Mount"
# ---------------------------------------------------------------
$dash = [char]8211
$titleText = "Book 3 " + $dash + " 50 " + $dash + " On The Olive Mount"
$foundTitle = $d.Content.Find.Execute($titleText, $true, $false, $false, $false, $false, $true, 1, $false, $titleText, 2)
Write-Output "Title replace: $foundTitle"

# ---------------------------------------------------------------
# Part 2: "schlimmer" paragraph
#   Merge "The German translated here as " + "schlimmer " + (start of
#   3rd run) into one run, then insert "It could be a temporary home. "
#   as its own run before "Combined with the friendly...".
# ---------------------------------------------------------------
$mergeOld = "schlimmer can also be"
$found2 = $d.Content.Find.Execute($mergeOld, $true, $false, $false, $false, $false, $true, 1, $false, $mergeOld, 2)
Write-Output "Schlimmer merge: $found2"

$rCombined = $d.Content
$rCombined.Find.Execute("Combined with the friendly") | Out-Null
$combinedStart = $rCombined.Start
$newRun2Text = "It could be a temporary home. "
$insPoint2 = $d.Range($combinedStart, $combinedStart)
$insPoint2.InsertBefore($newRun2Text)
$insRange2 = $d.Range($combinedStart, $combinedStart + $newRun2Text.Length)
$insRange2.Bold = 1
$insRange2.Bold = 0
Write-Output "Inserted 'It could be a temporary home.'"

# ---------------------------------------------------------------
# Part 3: "lick wounds" paragraph - append a new run at paragraph end
# ---------------------------------------------------------------
$rLick = $d.Content
$rLick.Find.Execute("lick wounds, used by David and Jesus. ") | Out-Null
$lickEnd = $rLick.End
$newRun3Text = "Zarathustra runs from what his  home, why is he so readily willing to abandon it? Because it is temporary, it is not his real home the cave. Perhaps one can imagine a student-dorm."
$insPoint3 = $d.Range($lickEnd, $lickEnd)
$insPoint3.InsertBefore($newRun3Text)
$insRange3 = $d.Range($lickEnd, $lickEnd + $newRun3Text.Length)
$insRange3.Bold = 1
$insRange3.Bold = 0
Write-Output "Inserted 'Zarathustra runs from...' "

# ---------------------------------------------------------------
# Part 4: "stern guest" paragraph
#   Merge the run ending in the two <w:br/> line breaks with the run
#   that begins "There, in his place of sanctuary...", then append a
#   trailing space plus a brand new run "A contrast to his Ape...".
# ---------------------------------------------------------------
$vtab = [char]11
$mergeOld4 = "stern guest, and am still fond of him; because he cleareth my house of flies, and quieteth many little noises. " + $vtab + $vtab + "There, in his place of sanctuary, which is not his temporary home"
$found4 = $d.Content.Find.Execute($mergeOld4, $true, $false, $false, $false, $false, $true, 1, $false, $mergeOld4, 2)
Write-Output "Stern guest merge: $found4"

$rGrateful = $d.Content
$rGrateful.Find.Execute("to finally now being grateful for what they do.") | Out-Null
$gratefulEnd = $rGrateful.End
$insSpace = $d.Range($gratefulEnd, $gratefulEnd)
$insSpace.InsertAfter(" ")

$newRun4Text = "A contrast to his Ape in the next session who provides negative utility because no one one honored him."
$insPoint4 = $d.Range($gratefulEnd + 1, $gratefulEnd + 1)
$insPoint4.InsertBefore($newRun4Text)
$insRange4 = $d.Range($gratefulEnd + 1, $gratefulEnd + 1 + $newRun4Text.Length)
$insRange4.Bold = 1
$insRange4.Bold = 0
Write-Output "Inserted 'A contrast to his Ape...'"

# ---------------------------------------------------------------
# Part 5: "Him whom I love. Whomever..." paragraph
#   Append a new run with more commentary text, then add three new
#   paragraphs after it: empty, "Next line,", empty.
# ---------------------------------------------------------------
$lastTextPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$loveEnd = $lastTextPara.Range.End - 1
$newRun5Text = "we are told many people Zarathustra loves. Still here it is whomever. Whomever that is Zarathustra loves more under winter. Under academic hardship, under physical hardship, times when they must burn internally lest they be consumed by cold. Now, Zarathustra is in his Olive-mount and he can mock his enemies with almost hungry zeal that his temporary home has fallen under the cold domain of winter. "
$insPoint5 = $d.Range($loveEnd, $loveEnd)
$insPoint5.InsertBefore($newRun5Text)
$insRange5 = $d.Range($loveEnd, $loveEnd + $newRun5Text.Length)
$insRange5.Bold = 1
$insRange5.Bold = 0
Write-Output "Inserted 'we are told many people...'"

$afterNewRun = $loveEnd + $newRun5Text.Length
$p1 = $d.Range($afterNewRun, $afterNewRun)
$p1.InsertParagraphAfter()
$p2 = $d.Range($afterNewRun, $afterNewRun)
$p2.InsertParagraphAfter()
$p3 = $d.Range($afterNewRun, $afterNewRun)
$p3.InsertParagraphAfter()

$nextLinePara = $d.Paragraphs.Item(91)
$nextLinePara.Range.Text = "Next line,"
Write-Output "Inserted trailing empty/Next line,/empty paragraphs"

Write-Output "All edits applied."
